$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price column (D) values are stored as text in the source data (scraped
# crypto prices, some containing multiple "." as thousands separators). Cells
# whose new value would otherwise be auto-recognised by Excel as a plain number
# are entered with a leading apostrophe (quote-prefix) so they stay text, exactly
# like the rest of the column.

# Row 2
$ws.Range("D2").Value = "25.699.55"
$ws.Range("E2").Value = "  -0.87%  "

# Row 3
$ws.Range("D3").Value = "1.629.98"
$ws.Range("E3").Value = "  -1.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'214.19"
$ws.Range("E5").Value = "  -1.00%  "

# Row 6
$ws.Range("E6").Value = "  -1.04%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("E8").Value = "  -0.83%  "

# Row 9
$ws.Range("D9").Value = "'0.0635"
$ws.Range("E9").Value = "  -1.53%  "

# Row 10
$ws.Range("D10").Value = "'19.49"
$ws.Range("E10").Value = "  -5.10%  "

# Row 11
$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "  -0.01%  "

# Row 12
$ws.Range("E12").Value = "  -0.93%  "

# Row 13
$ws.Range("D13").Value = "1.625.82"
$ws.Range("E13").Value = "  -1.48%  "

# Row 14
$ws.Range("D14").Value = "1.854.18"
$ws.Range("E14").Value = "  -1.11%  "

# Row 15
$ws.Range("E15").Value = "  -1.83%  "

# Row 16
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -1.36%  "

# Row 17
$ws.Range("D17").Value = "'63.11"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
$ws.Range("D18").Value = "25.722.44"

# Row 20
$ws.Range("E20").Value = "  +1.18%  "

# Row 21
$ws.Range("E21").Value = "  -0.16%  "

# Row 22
$ws.Range("D22").Value = "'9.93"
$ws.Range("E22").Value = "  -0.26%  "

# Row 23
$ws.Range("E23").Value = "  +0.97%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("E25").Value = "  -1.50%  "

# Row 26
$ws.Range("D26").Value = "'140.35"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27
$ws.Range("E27").Value = "  -3.50%  "

# Row 28
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("E29").Value = "  -0.36%  "

# Row 30
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "  -1.36%  "

# Row 31
$ws.Range("E31").Value = "  -2.69%  "

# Row 32
$ws.Range("E32").Value = "  +0.23%  "

# Row 33
$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("E36").Value = "  -1.53%  "

# Row 37
$ws.Range("E37").Value = "  -0.47%  "

# Row 38
$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "  -2.26%  "

# Row 39
$ws.Range("D39").Value = "1.103.50"
$ws.Range("E39").Value = "  -2.51%  "

# Row 40
$ws.Range("E40").Value = "  -1.22%  "

# Row 41
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("E42").Value = "  +0.76%  "

# Row 43
$ws.Range("D43").Value = "'99.62"
$ws.Range("E43").Value = "  +0.70%  "

# Row 44
$ws.Range("E44").Value = "  -1.16%  "

# Row 45
$ws.Range("D45").Value = "1.762.48"
$ws.Range("E45").Value = "  -1.22%  "

# Row 46
$ws.Range("E46").Value = "  -4.79%  "

# Row 47
$ws.Range("D47").Value = "'54.99"
$ws.Range("E47").Value = "  -1.56%  "

# Row 48
$ws.Range("E48").Value = "  -2.79%  "

# Row 49
$ws.Range("E49").Value = "  +4.57%  "

# Row 50 and 51: coin order swap (Cronos <-> EnergySwap) plus value updates
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.65"
$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0502"
$ws.Range("E51").Value = "  -0.58%  "
